$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header row's grouping merges have to go first - while they're still
# merged, only the anchor cell (H1/M1/Q1) of each merged block accepts a
# new value, so unmerge before filling in the rest of row 1.
$ws.Range("H1:L1").UnMerge()
$ws.Range("M1:P1").UnMerge()
$ws.Range("Q1:S1").UnMerge()

# Flatten the two-row (multi-index) header into a single header row.
# Row 1 previously held only the "group" labels (Tackles/Challenges/Blocks)
# over merged cells; now every column gets its own explicit label (mostly
# the same names that used to live in row 2), and row 2/3 become hidden
# detail/spacer rows instead.
$ws.Cells.Item(1, 1).Value = "Player ID"
$ws.Cells.Item(1, 2).Value = "Player"
$ws.Cells.Item(1, 3).Value = "#"
$ws.Cells.Item(1, 4).Value = "Nation"
$ws.Cells.Item(1, 5).Value = "Pos"
$ws.Cells.Item(1, 6).Value = "Age"
$ws.Cells.Item(1, 7).Value = "90s"
$ws.Cells.Item(1, 8).Value = "Tkl"
$ws.Cells.Item(1, 9).Value = "TklW"
$ws.Cells.Item(1, 10).Value = "Def 3rd"
$ws.Cells.Item(1, 11).Value = "Mid 3rd"
$ws.Cells.Item(1, 12).Value = "Att 3rd"
$ws.Cells.Item(1, 13).Value = "Cha"
$ws.Cells.Item(1, 14).Value = "Att"
$ws.Cells.Item(1, 15).Value = "Tkl%"
$ws.Cells.Item(1, 16).Value = "Lost"
$ws.Cells.Item(1, 17).Value = "Blocks"
$ws.Cells.Item(1, 18).Value = "Sh"
$ws.Cells.Item(1, 19).Value = "Pass"
$ws.Cells.Item(1, 20).Value = "Int"
$ws.Cells.Item(1, 21).Value = "Tkl+Int"
$ws.Cells.Item(1, 22).Value = "Clr"
$ws.Cells.Item(1, 23).Value = "Err"

# Old sub-header row, the spacer row beneath it, and the trailing totals
# row are kept for reference but hidden from view.
$ws.Rows.Item(2).Hidden = $true
$ws.Rows.Item(3).Hidden = $true
$ws.Rows.Item(15).Hidden = $true

# Fill in the previously-blank Tkl% cells (no tackle attempts recorded)
# with an explicit 0 instead of leaving them empty.
$ws.Cells.Item(5, 15).Value = 0
$ws.Cells.Item(11, 15).Value = 0
$ws.Cells.Item(12, 15).Value = 0
$ws.Cells.Item(13, 15).Value = 0
$ws.Cells.Item(14, 15).Value = 0

# Leave the selection where the editor last left it.
$ws.Range("O16").Select()
